$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "23.363.54"
$ws.Range("E2").Value = "  -0.78%  "

$ws.Range("D3").Value = "1.626.37"
$ws.Range("E3").Value = "  -0.91%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.000"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.04%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.9996"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.02%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "304.23"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.17%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3789"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.19%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "52.05"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.79%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.3630"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.59%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.231"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -4.12%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.08103"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.36%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.000"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.03%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "22.70"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -2.41%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.552"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.96%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.00001248"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -3.38%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "7.220"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -3.61%  "

$ws.Range("D17").Value = "1.628.55"
$ws.Range("E17").Value = "  -0.80%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "93.59"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.54%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06898"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.87%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.88"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -3.13%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.9997"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.05%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.421"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.77%  "

$ws.Range("D23").Value = "23.370.91"
$ws.Range("E23").Value = "  -0.79%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "12.73"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -2.15%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.254"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +3.03%  "

$ws.Range("E26").Value = "  +1.21%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "21.09"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.77%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "149.40"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.59%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.283"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.24%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "134.19"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.76%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.302"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -5.29%  "

$ws.Range("D32").Value = "1.806.46"
$ws.Range("E32").Value = "  -0.63%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.785"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.06%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "10.98"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +5.03%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.9515"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -2.77%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.02787"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.13%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.2524"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.95%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.08831"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.42%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "6.107"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -2.15%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.07144"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -4.74%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.360"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -3.03%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.7067"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.84%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "16.28"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.24%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "12.29"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -3.14%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.6453"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -3.10%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.322"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -2.10%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.9989"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.08%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.995"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.29%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.07995"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.01%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.200"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.65%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "125.77"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -4.38%  "
